$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.461.94"
$ws.Range("D3").Value = "2.475.90"
$ws.Range("E3").Value = "  +7.43%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.67"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +11.74%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.591"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.20%  "
$ws.Range("D9").Value = "2.474.66"
$ws.Range("E9").Value = "  +7.45%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.107"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.09%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.78"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.84%  "
$ws.Range("E12").Value = "  +1.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.355"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +7.66%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.56"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +15.01%  "
$ws.Range("D15").Value = "2.914.53"
$ws.Range("E15").Value = "  +7.48%  "
$ws.Range("D16").Value = "63.288.74"
$ws.Range("E16").Value = "  +5.91%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000145"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +10.31%  "
$ws.Range("D18").Value = "2.465.68"
$ws.Range("E18").Value = "  +7.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +9.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "344.81"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +11.35%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.32"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.88"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.58%  "
$ws.Range("E23").Value = "  +0.29%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.80"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.175"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.28%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("E27").Value = "  +13.89%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.18%  "
$ws.Range("E29").Value = "  +12.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.88"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +19.27%  "
$ws.Range("D31").Value = "0.0₃0817"
$ws.Range("E31").Value = "  +14.82%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.84"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +8.12%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "175.04"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.57%  "
$ws.Range("E34").Value = "  +12.78%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.401"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.39%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.93"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.19%  "
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.51"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +12.73%  "
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "369.90"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +19.27%  "
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.08%  "
$ws.Range("E41").Value = "  +14.85%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "40.26"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "151.98"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +12.17%  "
$ws.Range("E44").Value = "  +10.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.80"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +12.70%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.601"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.44%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0967"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.58%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0526"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.94%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.39"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +10.95%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0232"
$ws.Range("E50").Value = "  +3.56%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0227"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.27%  "
